$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Append the new "FALSO" group rows below the existing data (rows 50-65
#    first, column by column, then the lone "FALSO1" row at 49 - this
#    mirrors the order the labels were originally typed in).
for ($i = 2; $i -le 17; $i++) {
    $row = 48 + $i
    $ws.Cells.Item($row, 1).Value = "FALSO$i"
}

for ($i = 2; $i -le 17; $i++) {
    $row = 48 + $i
    $ws.Cells.Item($row, 2).Value = " "
}

$ws.Cells.Item(49, 1).Value = "FALSO1"

# 2) Update the label in B46: append a <br> line-break marker to the
#    "Cuenta de asign. de la renta primaria" entry (new_names column).
$ws.Range("B46").Value = "Cuenta de asign. <br> de la renta primaria"

# 3) Fill in the group counts (column C) for all of the new rows.
for ($i = 1; $i -le 17; $i++) {
    $row = 48 + $i
    $ws.Cells.Item($row, 3).Value = 0
}

# 4) Update the view's active selection to match the authored state.
$ws.Range("B46").Select()
